$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current ("before") values for rows 2 and 3
$B2 = $ws.Range("B2").Value()
$C2 = $ws.Range("C2").Value()
$E2 = $ws.Range("E2").Value()
$F2 = $ws.Range("F2").Value()

$B3 = $ws.Range("B3").Value()
$C3 = $ws.Range("C3").Value()
$E3 = $ws.Range("E3").Value()
$F3 = $ws.Range("F3").Value()

# --- Row 2: new "Arquivo" name, other columns come from old row 3 ---
$ws.Range("A2").Value = "copy_processo_1"
$ws.Range("B2").Value = $B3
$ws.Range("C2").Value = $C3
$ws.Range("E2").Value = $E3
# F2 becomes a date-looking text string; force text so Excel does not
# reinterpret it as a date serial, then restore the original cell style.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = $F3
$ws.Range("D2").Copy()
$ws.Range("F2").PasteSpecial(-4122)

# --- Row 3: new "Arquivo" name, other columns come from old row 2 ---
$ws.Range("A3").Value = "copy_processo_2"
$ws.Range("B3").Value = $B2
$ws.Range("C3").Value = $C2
# E3 becomes the text string "12723" rather than a number.
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "12723"
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = $F2
$ws.Range("D3").Copy()
$ws.Range("F3").PasteSpecial(-4122)

$excel.CutCopyMode = 0
